# Commentless GitLatch Commit @ 2024-6-25-13-51-28-418
#
# The reference diff for this commit touches nothing in the document's
# visible content (word/document.xml's <w:body> is byte-for-byte
# identical before/after -- the unified diff's only body-area hunk has
# no +/- lines, just context). Every real delta lives in package
# metadata that Word itself owns and that is not exposed anywhere on
# the Word object model:
#
#   * word/document.xml & word/styles.xml: the root element's
#     xmlns:w16du declaration and the "w16du" token in mc:Ignorable
#     are dropped.
#   * word/webextensions/taskpanes.xml: the persisted task-pane width
#     snaps from 438 back to the Office default of 350.
#   * word/webextensions/webextension1.xml: the webextension's id GUID
#     is reissued.
#
# These are all written by Word's internal add-in/task-pane host when
# it re-resolves a (here: unavailable/dev-registry) web add-in on
# open/save -- there is no Application/Document/TaskPane/WebExtension
# automation member for any of it (TaskPane exposes only Visible; the
# `wetp:taskpanes` / `we:webextension` parts aren't reachable through
# CustomXMLParts, WordOpenXML, or any other documented COM surface).
# A real macro cannot author these values, so nothing here is
# scriptable through legitimate COM automation.
#
# Touch the document the way a script driving this resave plausibly
# would (load -> no content mutation -> save back in place) without
# fabricating a body edit that isn't in the diff.
$d = $word.ActiveDocument

# Best-effort: if a given runtime *does* wire the task pane / web
# extension objects up to their persisted XML, pick up the intended
# values; on the stock Word OM (and this host) these members are
# absent/no-ops, so guard every call and continue quietly otherwise.
try {
    $taskPanes = $word.TaskPanes
    if ($taskPanes -and $taskPanes.Count -gt 0) {
        for ($i = 1; $i -le $taskPanes.Count; $i++) {
            $pane = $taskPanes.Item($i)
            if ($pane) {
                try { $pane.Width = 350 } catch { }
            }
        }
    }
} catch { }

try {
    $addIns = $word.COMAddIns
    if ($addIns -and $addIns.Count -gt 0) {
        for ($i = 1; $i -le $addIns.Count; $i++) {
            $addIn = $addIns.Item($i)
            if ($addIn) {
                try { $addIn.Connect = $addIn.Connect } catch { }
            }
        }
    }
} catch { }

# No in-body content changed in the source revision, so there is
# nothing to Find/Replace. Persist the (unchanged) document back,
# mirroring the no-op content save that produced the captured diff.
$d.Save()
